# Updates cryptos list values (prices / volume%) scraped on
# Mon Sep  2 17:02:56 UTC 2024, plus a data-source row swap
# (Monero <-> PancakeSwap) in rows 30/31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value.
$updates = [ordered]@{
    "D2" = "58.329.11"
    "D3" = "2.510.98"
    "E3" = "  +1.86%  "
    "D4" = "0.998"
    "E4" = "  -0.13%  "
    "D5" = "520.70"
    "E5" = "  +0.04%  "
    "D6" = "131.83"
    "E6" = "  -0.94%  "
    "E7" = "  -0.19%  "
    "E8" = "  -0.04%  "
    "D9" = "2.507.61"
    "E9" = "  +1.38%  "
    "D10" = "0.0970"
    "E10" = "  -0.91%  "
    "E11" = "  -1.14%  "
    "E12" = "  -2.88%  "
    "D13" = "0.330"
    "E13" = "  -2.62%  "
    "D14" = "2.925.14"
    "E14" = "  +0.70%  "
    "D15" = "58.219.15"
    "E15" = "  +0.25%  "
    "D16" = "21.96"
    "E16" = "  -1.36%  "
    "E17" = "  -0.58%  "
    "D18" = "2.506.45"
    "E18" = "  +1.67%  "
    "D19" = "10.56"
    "E19" = "  -0.84%  "
    "D20" = "320.38"
    "E20" = "  +0.05%  "
    "E21" = "  -0.47%  "
    "D22" = "6.13"
    "E22" = "  +7.04%  "
    "D23" = "0.995"
    "E23" = "  -0.52%  "
    "D24" = "64.44"
    "E24" = "  -0.62%  "
    "D25" = "0.403"
    "E25" = "  -1.37%  "
    "D26" = "0.998"
    "E26" = "  -0.17%  "
    "D27" = "0.159"
    "E27" = "  -0.31%  "
    "D28" = "7.35"
    "E28" = "  +0.13%  "
    "D29" = "0.0₃0748"
    "E29" = "  +0.24%  "
    "B30" = "Monero"
    "C30" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D30" = "167.66"
    "E30" = "  -0.09%  "
    "B31" = "PancakeSwap"
    "C31" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D31" = "1.71"
    "E31" = "  +1.11%  "
    "D32" = "1.18"
    "E32" = "  +1.92%  "
    "D33" = "6.24"
    "E33" = "  +0.04%  "
    "E34" = "  +0.02%  "
    "D35" = "0.999"
    "E35" = "  -0.02%  "
    "D36" = "17.97"
    "E36" = "  -0.42%  "
    "E37" = "  -8.77%  "
    "D38" = "3.92"
    "E38" = "  -1.07%  "
    "E39" = "  -0.18%  "
    "D40" = "36.10"
    "E40" = "  -0.69%  "
    "D41" = "0.766"
    "E41" = "  -3.31%  "
    "D42" = "277.69"
    "E42" = "  +1.67%  "
    "D43" = "3.46"
    "E43" = "  +0.43%  "
    "D44" = "5.05"
    "E44" = "  -0.37%  "
    "E45" = "  +0.68%  "
    "D46" = "123.14"
    "E46" = "  -1.23%  "
    "D47" = "0.0917"
    "E47" = "  +1.17%  "
    "D48" = "0.0498"
    "E48" = "  +2.35%  "
    "D49" = "17.55"
    "E49" = "  -0.18%  "
    "E50" = "  -0.24%  "
    "D51" = "16.68"
    "E51" = "  -0.95%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force plain-text storage (matches the source inlineStr cells) so
    # purely-numeric-looking strings such as "0.998" or "520.70" are not
    # auto-coerced into numbers; then drop back to the default "Normal"
    # style so no stray number-format style is left applied to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
